$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the "Tv20 C4" row (row 2), shifting the rows below it up.
$ws.Rows(2).Delete()

# Leave the selection on B3, matching where the user ended up after the edit.
$ws.Range("B3").Select()
